# Auto-generated script applying market-price/profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 88
$ws.Range("I4").Value = 88
$ws.Range("K4").Value = 88
$ws.Range("M4").Value = 26
$ws.Range("H32").Value = 11870.5
$ws.Range("I32").Value = 12255.25
$ws.Range("J32").Value = 11101
$ws.Range("K32").Value = 12255.25
$ws.Range("L32").Value = 11101
$ws.Range("M32").Value = -11929.25
$ws.Range("N32").Value = -11753
$ws.Range("H70").Value = 6198.2
$ws.Range("I70").Value = 1194.5
$ws.Range("K70").Value = 3583.5
$ws.Range("M70").Value = -3313.5
$ws.Range("H73").Value = 6198.2
$ws.Range("I73").Value = 1194.5
$ws.Range("K73").Value = 3583.5
$ws.Range("M73").Value = -2647.5
$ws.Range("H82").Value = 1699.75
$ws.Range("I82").Value = 1699.75
$ws.Range("K82").Value = 5099.25
$ws.Range("M82").Value = -4693.25
$ws.Range("H85").Value = 1699.75
$ws.Range("I85").Value = 1699.75
$ws.Range("K85").Value = 5099.25
$ws.Range("M85").Value = -3695.25
$ws.Range("H103").Value = 3847.5
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 3847.5
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 11542.5
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -12714.5
$ws.Range("H116").Value = 16001
$ws.Range("I116").Value = 14601
$ws.Range("K116").Value = 14601
$ws.Range("M116").Value = -11159
$ws.Range("H137").Value = 3198.0286
$ws.Range("I137").Value = 2831
$ws.Range("K137").Value = 8493
$ws.Range("M137").Value = -5943

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 230.57143
$ws.Range("I5").Value = 172.875
$ws.Range("J5").Value = 307.5
$ws.Range("K5").Value = 172.875
$ws.Range("L5").Value = 307.5
$ws.Range("M5").Value = -60.875
$ws.Range("N5").Value = -531.5
$ws.Range("H32").Value = 4058.625
$ws.Range("I32").Value = 2973.6086
$ws.Range("K32").Value = 2973.6086
$ws.Range("M32").Value = -2686.6086
$ws.Range("H61").Value = 3642.762
$ws.Range("I61").Value = 3802.5945
$ws.Range("K61").Value = 3802.5945
$ws.Range("M61").Value = -3590.5945
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H112").Value = 62916.668
$ws.Range("J112").Value = 62916.668
$ws.Range("L112").Value = 62916.668
$ws.Range("N112").Value = -65870.66800000001
$ws.Range("H122").Value = 2888.6667
$ws.Range("I122").Value = 1706.25
$ws.Range("K122").Value = 5118.75
$ws.Range("M122").Value = -2668.75
$ws.Range("H136").Value = 3642.762
$ws.Range("I136").Value = 3802.5945
$ws.Range("K136").Value = 11407.7835
$ws.Range("M136").Value = -8857.783500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 230.57143
$ws.Range("I4").Value = 172.875
$ws.Range("J4").Value = 307.5
$ws.Range("K4").Value = 172.875
$ws.Range("L4").Value = 307.5
$ws.Range("M4").Value = -57.875
$ws.Range("N4").Value = -537.5
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -11058
$ws.Range("H94").Value = 641.9231
$ws.Range("I94").Value = 654.6
$ws.Range("J94").Value = 599.6667
$ws.Range("K94").Value = 654.6
$ws.Range("L94").Value = 599.6667
$ws.Range("M94").Value = -203.6
$ws.Range("N94").Value = -1501.6667
$ws.Range("H134").Value = 1955.0312
$ws.Range("I134").Value = 1108.3214
$ws.Range("J134").Value = 7882
$ws.Range("K134").Value = 3324.9642
$ws.Range("L134").Value = 23646
$ws.Range("M134").Value = -789.9642000000003
$ws.Range("N134").Value = -28716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 705.2
$ws.Range("I7").Value = 718.75
$ws.Range("K7").Value = 718.75
$ws.Range("M7").Value = -605.75
$ws.Range("H22").Value = 3232.5715
$ws.Range("J22").Value = 4365.6
$ws.Range("L22").Value = 4365.6
$ws.Range("N22").Value = -5065.6
$ws.Range("H58").Value = 3277.8333
$ws.Range("I58").Value = 1552.44
$ws.Range("J58").Value = 11904.8
$ws.Range("K58").Value = 1552.44
$ws.Range("L58").Value = 11904.8
$ws.Range("M58").Value = -1349.44
$ws.Range("N58").Value = -12310.8
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H130").Value = 41852.668
$ws.Range("J130").Value = 50000
$ws.Range("L130").Value = 50000
$ws.Range("N130").Value = -60040
$ws.Range("H132").Value = 1768.425
$ws.Range("I132").Value = 1223
$ws.Range("J132").Value = 4339.7144
$ws.Range("K132").Value = 3669
$ws.Range("L132").Value = 13019.1432
$ws.Range("M132").Value = -1139
$ws.Range("N132").Value = -18079.1432
$ws.Range("H136").Value = 3277.8333
$ws.Range("I136").Value = 1552.44
$ws.Range("J136").Value = 11904.8
$ws.Range("K136").Value = 4657.32
$ws.Range("L136").Value = 35714.39999999999
$ws.Range("M136").Value = -2107.32
$ws.Range("N136").Value = -40814.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H86").Value = 2400
$ws.Range("J86").Value = 3474.6667
$ws.Range("L86").Value = 10424.0001
$ws.Range("N86").Value = -12796.0001
$ws.Range("H89").Value = 2400
$ws.Range("J89").Value = 3474.6667
$ws.Range("L89").Value = 31272.0003
$ws.Range("N89").Value = -43128.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2527.7778
$ws.Range("I43").Value = 2527.7778
$ws.Range("K43").Value = 2527.7778
$ws.Range("M43").Value = -2376.7778
$ws.Range("H99").Value = 34170.215
$ws.Range("I99").Value = 33624.145
$ws.Range("J99").Value = 34716.285
$ws.Range("K99").Value = 33624.145
$ws.Range("L99").Value = 34716.285
$ws.Range("M99").Value = -31378.145
$ws.Range("N99").Value = -39208.285
$ws.Range("H113").Value = 6012.2666
$ws.Range("I113").Value = 5953.1
$ws.Range("J113").Value = 6130.6
$ws.Range("K113").Value = 5953.1
$ws.Range("L113").Value = 6130.6
$ws.Range("M113").Value = -3783.1
$ws.Range("N113").Value = -10470.6
$ws.Range("H126").Value = 3263.0527
$ws.Range("I126").Value = 2387.2942
$ws.Range("K126").Value = 7161.882599999999
$ws.Range("M126").Value = -4691.882599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5738.55
$ws.Range("J22").Value = 7538.4287
$ws.Range("L22").Value = 7538.4287
$ws.Range("N22").Value = -8128.4287
$ws.Range("H27").Value = 5738.55
$ws.Range("J27").Value = 7538.4287
$ws.Range("L27").Value = 7538.4287
$ws.Range("N27").Value = -7752.4287
$ws.Range("H122").Value = 4783.3335
$ws.Range("I122").Value = 3857.3076
$ws.Range("J122").Value = 6288.125
$ws.Range("K122").Value = 11571.9228
$ws.Range("L122").Value = 18864.375
$ws.Range("M122").Value = -9121.9228
$ws.Range("N122").Value = -23764.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 50615
$ws.Range("J105").Value = 50615
$ws.Range("L105").Value = 50615
$ws.Range("N105").Value = -57603
$ws.Range("H107").Value = 1278.3
$ws.Range("I107").Value = 989
$ws.Range("K107").Value = 2967
$ws.Range("M107").Value = -1047
$ws.Range("H108").Value = 78947
$ws.Range("J108").Value = 78947
$ws.Range("L108").Value = 78947
$ws.Range("N108").Value = -86627
